$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: convert PartNumber (C5) and Modelo (D5) from text to numeric values
$ws.Range("C5").Value = 53495349
$ws.Range("D5").Value = 5980
$ws.Range("E5").Value = 45960.63015388889
$ws.Range("E5").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Add new row 6 data for the new piece (LONGUERONE CPL LT)
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "LONGUERONE CPL LT"
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = "53489572"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "5980"

$ws.Range("E6").Value = 45969.89631894375
$ws.Range("E6").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("F6").Value = "dados/peca_5/txt"
$ws.Range("G6").Value = "Ativa"
